$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post that used to occupy row 847 ("恐れの中に生きる者は決して自由になれない")
# was removed. Delete that entire row; Excel will automatically shift all
# subsequent rows (848:883) up by one, renumbering them to 847:882 and
# shrinking the sheet's used range/dimension accordingly.
$ws.Rows(847).Delete()
